# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strike) values computed for rows 2-22 (column G)
$kValues = @(4, 3, 7, 3, 5, 6, 1, 3, 6, 1, 2, 8, 7, 2, 4, 1, 3, 0, 2, 1, 5)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
